$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("goods init")
$ws1.Range("A10").Value = "salt"
$ws1.Range("B10").Value = 1
$ws1.Range("C10").Value = "kg"

$ws2 = $wb.Worksheets.Item("product init")
$ws2.Range("A10").Value = "croissant"
$ws2.Range("B10").Value = 0.4
$ws2.Range("C10").Value = 0
